# Update (Removed Auto Arima)
# Updates forecast figures on the "Forecast Comparison" sheet (Prophet / Amazon
# Mean / P70 / P80 / P90 forecast columns, rows 2-17) and the corresponding
# aggregate figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Row => @(Prophet Forecast(C), Amazon Mean Forecast(D), Amazon P70 Forecast(E), Amazon P80 Forecast(F), Amazon P90 Forecast(G))
$forecastRows = @{
    2  = @(65, 46, 55, 65, 80)
    3  = @(57, 38, 46, 56, 71)
    4  = @(61, 35, 42, 51, 65)
    5  = @(72, 36, 44, 53, 68)
    6  = @(75, 36, 44, 55, 71)
    7  = @(64, 36, 44, 54, 70)
    8  = @(41, 36, 44, 55, 74)
    9  = @(22, 35, 43, 55, 74)
    10 = @(16, 35, 42, 53, 69)
    11 = @(22, 35, 43, 54, 71)
    12 = @(29, 35, 42, 54, 73)
    13 = @(33, 36, 44, 56, 75)
    14 = @(38, 34, 42, 53, 71)
    15 = @(52, 33, 41, 53, 74)
    16 = @(71, 33, 40, 51, 70)
    17 = @(87, 32, 39, 51, 70)
}

foreach ($row in $forecastRows.Keys) {
    $vals = $forecastRows[$row]
    $wsForecast.Cells.Item($row, 3).Value = $vals[0]  # C - Prophet Forecast
    $wsForecast.Cells.Item($row, 4).Value = $vals[1]  # D - Amazon Mean Forecast
    $wsForecast.Cells.Item($row, 5).Value = $vals[2]  # E - Amazon P70 Forecast
    $wsForecast.Cells.Item($row, 6).Value = $vals[3]  # F - Amazon P80 Forecast
    $wsForecast.Cells.Item($row, 7).Value = $vals[4]  # G - Amazon P90 Forecast
}

# Summary sheet aggregate values (stored as text, matching the original cell type)
$summaryRows = @{
    9  = "805"  # Total Forecast (16 Weeks)
    10 = "457"  # Total Forecast (8 Weeks)
    11 = "255"  # Total Forecast (4 Weeks)
    12 = "87"   # Max Forecast
    14 = "16"   # Min Forecast
}

foreach ($row in $summaryRows.Keys) {
    $cell = $wsSummary.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryRows[$row]
}
